# "Generate Report for Handback"
# Populates the "Latest Target File" (F) and "Latest Handback File" (G)
# columns for each locale sheet, marks the rows as handed back (Status,
# column C) and stamps the handback datetime (column H).

$wb = $excel.ActiveWorkbook

# The Overview sheet mirrors each locale's status in column B (zh-cn) and
# column C (de-de) via the same shared-string text, so it needs to be
# stamped with the new status as well.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime,
        [string]$MdUrl1,
        [string]$XlfUrl1,
        [string]$MdUrl2,
        [string]$XlfUrl2
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # --- Row 2 : 4c46dfff-af03-4564-b2a8-b45668200950 -------------------
    $ws.Range("C2").Value = "Handed back: in sync with en-US"

    $mdName1 = $ws.Range("A2").Value2
    $ws.Hyperlinks.Add($ws.Range("F2"), $MdUrl1, "", "", $mdName1) | Out-Null

    $xlfName1 = $ws.Range("D2").Value2
    $ws.Hyperlinks.Add($ws.Range("G2"), $XlfUrl1, "", "", $xlfName1) | Out-Null

    $ws.Range("H2").Value = $HandbackDateTime

    # --- Row 3 : fb7702b0-52c6-44c3-bce9-be768f8c0aae -------------------
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    $mdName2 = $ws.Range("A3").Value2
    $ws.Hyperlinks.Add($ws.Range("F3"), $MdUrl2, "", "", $mdName2) | Out-Null

    $xlfName2 = $ws.Range("D3").Value2
    $ws.Hyperlinks.Add($ws.Range("G3"), $XlfUrl2, "", "", $xlfName2) | Out-Null

    $ws.Range("H3").Value = $HandbackDateTime
}

Update-LocaleSheet "zh-cn" "2016-03-21 04:39:41" `
    "https://github.com/OpenLocalizationTest/oltest/blob/63315741237e9267de526a187d800dde99e0efb4/e2e/4c46dfff-af03-4564-b2a8-b45668200950.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2155098fffe5ec0a80b025c0714103efff50a12/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/4c46dfff-af03-4564-b2a8-b45668200950.47a0ef89d373fe19bcdd383c2d85ee9dfbfa71e4.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/63315741237e9267de526a187d800dde99e0efb4/e2e/fb7702b0-52c6-44c3-bce9-be768f8c0aae.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2155098fffe5ec0a80b025c0714103efff50a12/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/fb7702b0-52c6-44c3-bce9-be768f8c0aae.03e2e94da20a457dc856a1ffefc5e3abd1fe94eb.zh-cn.xlf"

Update-LocaleSheet "de-de" "2016-03-21 04:39:56" `
    "https://github.com/OpenLocalizationTest/oltest/blob/63315741237e9267de526a187d800dde99e0efb4/e2e/4c46dfff-af03-4564-b2a8-b45668200950.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e826aa3995aa0347ce07e1da6697546b9677613/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/4c46dfff-af03-4564-b2a8-b45668200950.47a0ef89d373fe19bcdd383c2d85ee9dfbfa71e4.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/63315741237e9267de526a187d800dde99e0efb4/e2e/fb7702b0-52c6-44c3-bce9-be768f8c0aae.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e826aa3995aa0347ce07e1da6697546b9677613/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/fb7702b0-52c6-44c3-bce9-be768f8c0aae.03e2e94da20a457dc856a1ffefc5e3abd1fe94eb.de-de.xlf"

Write-Host "Handback report generated."
